$d = $word.ActiveDocument

# Locate the paragraph that contains the text
# "...Teste qui-quadrado (nominal ou ordinal)" so we can reliably find the
# empty paragraph that immediately follows it, regardless of paragraph index.
$searchRange = $d.Content
$searchRange.Find.ClearFormatting()
$found = $searchRange.Find.Execute("qui-quadrado (nominal ou ordinal)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor text for the target paragraph."
}

$anchorPara = $searchRange.Paragraphs(1)

# The paragraph right after the anchor is the empty paragraph that needs the
# new left indent (0.25in = 18pt = 360 twips).
$targetPara = $anchorPara.Next()
$targetPara.LeftIndent = 18

# Duplicate that empty paragraph: insert a new paragraph right after it that
# carries the same (now-updated) paragraph formatting.
$targetPara.Range.InsertParagraphAfter()
